$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-25 with the refreshed candidate pipeline data
$ws.Cells.Item(2, 1).Value = 196
$ws.Cells.Item(2, 2).Value = "Blockaid"
$ws.Cells.Item(2, 3).Value = "Enterprise Account Executive (Fintech)"
$ws.Cells.Item(2, 4).Value = "Connor Yakushi"
$ws.Cells.Item(2, 5).Value = "CV Sent"

$ws.Cells.Item(3, 1).Value = 196
$ws.Cells.Item(3, 2).Value = "Blockaid"
$ws.Cells.Item(3, 3).Value = "Enterprise Account Executive (Fintech)"
$ws.Cells.Item(3, 4).Value = "Madyson Almeida"
$ws.Cells.Item(3, 5).Value = "2nd Interview"

$ws.Cells.Item(4, 1).Value = 196
$ws.Cells.Item(4, 2).Value = "Blockaid"
$ws.Cells.Item(4, 3).Value = "Enterprise Account Executive (Fintech)"
$ws.Cells.Item(4, 4).Value = "Mike Gomez"
$ws.Cells.Item(4, 5).Value = "1st Interview"

$ws.Cells.Item(5, 1).Value = 484
$ws.Cells.Item(5, 2).Value = "Cognition AI"
$ws.Cells.Item(5, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(5, 4).Value = "Andrew Rydalch"
$ws.Cells.Item(5, 5).Value = "1st Interview"

$ws.Cells.Item(6, 1).Value = 484
$ws.Cells.Item(6, 2).Value = "Cognition AI"
$ws.Cells.Item(6, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(6, 4).Value = "Mikaela Stamas"
$ws.Cells.Item(6, 5).Value = "4th Interview"

$ws.Cells.Item(7, 1).Value = 484
$ws.Cells.Item(7, 2).Value = "Cognition AI"
$ws.Cells.Item(7, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(7, 4).Value = "Tiffany Shih"
$ws.Cells.Item(7, 5).Value = "1st Interview"

$ws.Cells.Item(8, 1).Value = 591
$ws.Cells.Item(8, 2).Value = "Doxel.ai"
$ws.Cells.Item(8, 3).Value = "CS1 Doxel - Enterprise AE Northeast U.S"
$ws.Cells.Item(8, 4).Value = "Mike Gomez"
$ws.Cells.Item(8, 5).Value = "1st Interview"

$ws.Cells.Item(9, 1).Value = 633
$ws.Cells.Item(9, 2).Value = "Factory"
$ws.Cells.Item(9, 3).Value = "CS1 Factory - Enterprise AE"
$ws.Cells.Item(9, 4).Value = "Tiffany Shih"
$ws.Cells.Item(9, 5).Value = "1st Interview"

$ws.Cells.Item(10, 1).Value = 652
$ws.Cells.Item(10, 2).Value = "Dash0"
$ws.Cells.Item(10, 3).Value = "Enterprise AE PST"
$ws.Cells.Item(10, 4).Value = "Mike Gomez"
$ws.Cells.Item(10, 5).Value = "CV Sent"

$ws.Cells.Item(11, 1).Value = 750
$ws.Cells.Item(11, 2).Value = "Novee.io"
$ws.Cells.Item(11, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(11, 4).Value = "Garrett Kingston"
$ws.Cells.Item(11, 5).Value = "CV Sent"

$ws.Cells.Item(12, 1).Value = 750
$ws.Cells.Item(12, 2).Value = "Novee.io"
$ws.Cells.Item(12, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(12, 4).Value = "Tricia Rupp"
$ws.Cells.Item(12, 5).Value = "2nd Interview"

$ws.Cells.Item(13, 1).Value = 750
$ws.Cells.Item(13, 2).Value = "Novee.io"
$ws.Cells.Item(13, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(13, 4).Value = "Danny Nia"
$ws.Cells.Item(13, 5).Value = "CV Sent"

$ws.Cells.Item(14, 1).Value = 757
$ws.Cells.Item(14, 2).Value = "Blockaid"
$ws.Cells.Item(14, 3).Value = "SDR (Singapore)"
$ws.Cells.Item(14, 4).Value = "Annisa Mareizky"
$ws.Cells.Item(14, 5).Value = "1st Interview"

$ws.Cells.Item(15, 1).Value = 757
$ws.Cells.Item(15, 2).Value = "Blockaid"
$ws.Cells.Item(15, 3).Value = "SDR (Singapore)"
$ws.Cells.Item(15, 4).Value = "Dinie Mifdhal"
$ws.Cells.Item(15, 5).Value = "4th Interview"

$ws.Cells.Item(16, 1).Value = 757
$ws.Cells.Item(16, 2).Value = "Blockaid"
$ws.Cells.Item(16, 3).Value = "SDR (Singapore)"
$ws.Cells.Item(16, 4).Value = "Dinie Mifdhal"
$ws.Cells.Item(16, 5).Value = "4th Interview"

$ws.Cells.Item(17, 1).Value = 757
$ws.Cells.Item(17, 2).Value = "Blockaid"
$ws.Cells.Item(17, 3).Value = "SDR (Singapore)"
$ws.Cells.Item(17, 4).Value = "Hayoung Kim"
$ws.Cells.Item(17, 5).Value = "1st Interview"

$ws.Cells.Item(18, 1).Value = 766
$ws.Cells.Item(18, 2).Value = "Cogent Security"
$ws.Cells.Item(18, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(18, 4).Value = "Max B."
$ws.Cells.Item(18, 5).Value = "CV Sent"

$ws.Cells.Item(19, 1).Value = 779
$ws.Cells.Item(19, 2).Value = "Energy Robotics"
$ws.Cells.Item(19, 3).Value = "SDR London"
$ws.Cells.Item(19, 4).Value = "Fouad Abou-Steit"
$ws.Cells.Item(19, 5).Value = "2nd Interview"

$ws.Cells.Item(20, 1).Value = 783
$ws.Cells.Item(20, 2).Value = "Port"
$ws.Cells.Item(20, 3).Value = "Mid-Market AE"
$ws.Cells.Item(20, 4).Value = "Andy Evans"
$ws.Cells.Item(20, 5).Value = "1st Interview"

$ws.Cells.Item(21, 1).Value = 783
$ws.Cells.Item(21, 2).Value = "Port"
$ws.Cells.Item(21, 3).Value = "Mid-Market AE"
$ws.Cells.Item(21, 4).Value = "Shahz Shuja"
$ws.Cells.Item(21, 5).Value = "3rd Interview"

$ws.Cells.Item(22, 1).Value = 816
$ws.Cells.Item(22, 2).Value = "Allium"
$ws.Cells.Item(22, 3).Value = "Enterprise Account Executive"
$ws.Cells.Item(22, 4).Value = "Erik Hug"
$ws.Cells.Item(22, 5).Value = "1st Interview"

$ws.Cells.Item(23, 1).Value = 824
$ws.Cells.Item(23, 2).Value = "Blockaid"
$ws.Cells.Item(23, 3).Value = "Technical Account Manager"
$ws.Cells.Item(23, 4).Value = "Jeff Cooperstein"
$ws.Cells.Item(23, 5).Value = "CV Sent"

$ws.Cells.Item(24, 1).Value = 847
$ws.Cells.Item(24, 2).Value = "Simile.ai"
$ws.Cells.Item(24, 3).Value = "CS1 Simile.ai - Enterprise AE x2"
$ws.Cells.Item(24, 4).Value = "Amelia Silverwood"
$ws.Cells.Item(24, 5).Value = "CV Sent"

$ws.Cells.Item(25, 1).Value = 865
$ws.Cells.Item(25, 2).Value = "Accel Data"
$ws.Cells.Item(25, 3).Value = "Enterprise Account Executive"
$ws.Cells.Item(25, 4).Value = "Sultan Miah"
$ws.Cells.Item(25, 5).Value = "1st Interview"

# Remove now-unused trailing rows 26-29 (table shrank from 29 to 25 rows)
$ws.Range("A26:E29").ClearContents() | Out-Null
